$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style) from the last existing row (A301)
# down through the newly appended rows (A302:A328).
$ws.Range("A301").Copy()
$ws.Range("A302:A328").PasteSpecial(-4122)

# New data rows to append: row, date serial (col A), col B, col C, col D
$data = @(
    @(302, 44376, 0, 0, 0),
    @(303, 44377, 0, 0, 0),
    @(304, 44378, 0, 0, 0),
    @(305, 44379, 3, 3, 36.46086533787069),
    @(306, 44380, 0, 3, 36.46086533787069),
    @(307, 44381, 0, 3, 36.46086533787069),
    @(308, 44382, 0, 3, 36.46086533787069),
    @(309, 44383, 0, 3, 36.46086533787069),
    @(310, 44384, 0, 3, 36.46086533787069),
    @(311, 44385, 0, 3, 36.46086533787069),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 1, 1, 12.15362177929023),
    @(317, 44391, 0, 1, 12.15362177929023),
    @(318, 44392, 0, 1, 12.15362177929023),
    @(319, 44393, 0, 1, 12.15362177929023),
    @(320, 44394, 0, 1, 12.15362177929023),
    @(321, 44395, 0, 1, 12.15362177929023),
    @(322, 44396, 0, 1, 12.15362177929023),
    @(323, 44397, 0, 0, 0),
    @(324, 44398, 0, 0, 0),
    @(325, 44399, 0, 0, 0),
    @(326, 44400, 0, 0, 0),
    @(327, 44401, 0, 0, 0),
    @(328, 44402, 1, 1, 12.15362177929023)
)

foreach ($r in $data) {
    $rowIndex = $r[0]
    $ws.Cells.Item($rowIndex, 1).Value = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]
}
